$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.535309666666667
$ws.Range("H2").Value = 16.605929
$ws.Range("I2").Value = 0.1390876011186406
$ws.Range("J2").Value = 0.1461320332765693
$ws.Range("M2").Value = 7.207612333333334
$ws.Range("N2").Value = 21.622837
$ws.Range("O2").Value = 0.1890664495410015
$ws.Range("P2").Value = 0.2007919674510081
$ws.Range("Q2").Value = 39.89636622228589
$ws.Range("R2").Value = 359.067296000573
$ws.Range("S2").Value = 0.02629679891867642
$ws.Range("T2").Value = 0.02934213846921852

$ws.Range("G3").Value = 5.535309666666667
$ws.Range("H3").Value = 16.605929
$ws.Range("I3").Value = 0.1390876011186406
$ws.Range("J3").Value = 0.1461320332765693
$ws.Range("O3").Value = 0.2101379387118686
$ws.Range("P3").Value = 0.2231702676624555
$ws.Range("Q3").Value = 44.34282327931933
$ws.Range("R3").Value = 399.085409513874
$ws.Range("S3").Value = 0.02922758179944972
$ws.Range("T3").Value = 0.03261232498039082

$ws.Range("G4").Value = 5.535309666666667
$ws.Range("H4").Value = 16.605929
$ws.Range("I4").Value = 0.1390876011186406
$ws.Range("J4").Value = 0.1461320332765693
$ws.Range("M4").Value = 9.361017666666667
$ws.Range("N4").Value = 28.083053
$ws.Range("O4").Value = 0.2455534915692039
$ws.Range("P4").Value = 0.2607822213107805
$ws.Range("Q4").Value = 51.81613158013744
$ws.Range("R4").Value = 466.345184221237
$ws.Range("S4").Value = 0.03415344608866692
$ws.Range("T4").Value = 0.03810863624252463

$ws.Range("G5").Value = 5.535309666666667
$ws.Range("H5").Value = 16.605929
$ws.Range("I5").Value = 0.1390876011186406
$ws.Range("J5").Value = 0.1461320332765693
$ws.Range("M5").Value = 6.6785765
$ws.Range("N5").Value = 13.357153
$ws.Range("O5").Value = 0.1751890485290578
$ws.Range("P5").Value = 0.1240359454411156
$ws.Range("Q5").Value = 36.96798906002283
$ws.Range("R5").Value = 221.807934360137
$ws.Range("S5").Value = 0.02436662450216376
$ws.Range("T5").Value = 0.01812562490669183

$ws.Range("G6").Value = 5.535309666666667
$ws.Range("H6").Value = 16.605929
$ws.Range("I6").Value = 0.1390876011186406
$ws.Range("J6").Value = 0.1461320332765693
$ws.Range("M6").Value = 6.864003333333333
$ws.Range("N6").Value = 20.59201
$ws.Range("O6").Value = 0.1800530716488682
$ws.Range("P6").Value = 0.1912195981346403
$ws.Range("Q6").Value = 37.99438400303222
$ws.Range("R6").Value = 341.94945602729
$ws.Range("S6").Value = 0.0250431498096838
$ws.Range("T6").Value = 0.02794330867774347

$ws.Range("I7").Value = 0.2300921801028976
$ws.Range("J7").Value = 0.2417457620165159
$ws.Range("M7").Value = 7.207612333333334
$ws.Range("N7").Value = 21.622837
$ws.Range("O7").Value = 0.1890664495410015
$ws.Range("P7").Value = 0.2007919674510081
$ws.Range("Q7").Value = 66.00043288142578
$ws.Range("R7").Value = 594.0038959328321
$ws.Range("S7").Value = 0.04350271155920352
$ws.Range("T7").Value = 0.04854060717823941

$ws.Range("I8").Value = 0.2300921801028976
$ws.Range("J8").Value = 0.2417457620165159
$ws.Range("O8").Value = 0.2101379387118686
$ws.Range("P8").Value = 0.2231702676624555
$ws.Range("S8").Value = 0.04835109644054291
$ws.Range("T8").Value = 0.05395046641549012

$ws.Range("I9").Value = 0.2300921801028976
$ws.Range("J9").Value = 0.2417457620165159
$ws.Range("M9").Value = 9.361017666666667
$ws.Range("N9").Value = 28.083053
$ws.Range("O9").Value = 0.2455534915692039
$ws.Range("P9").Value = 0.2607822213107805
$ws.Range("Q9").Value = 85.7192631398009
$ws.Range("R9").Value = 771.473368258208
$ws.Range("S9").Value = 0.05649993820703661
$ws.Range("T9").Value = 0.06304299681113434

$ws.Range("I10").Value = 0.2300921801028976
$ws.Range("J10").Value = 0.2417457620165159
$ws.Range("M10").Value = 6.6785765
$ws.Range("N10").Value = 13.357153
$ws.Range("O10").Value = 0.1751890485290578
$ws.Range("P10").Value = 0.1240359454411156
$ws.Range("Q10").Value = 61.15602777263467
$ws.Range("R10").Value = 366.936166635808
$ws.Range("S10").Value = 0.04030963010620322
$ws.Range("T10").Value = 0.02998516414810147

$ws.Range("I11").Value = 0.2300921801028976
$ws.Range("J11").Value = 0.2417457620165159
$ws.Range("M11").Value = 6.864003333333333
$ws.Range("N11").Value = 20.59201
$ws.Range("O11").Value = 0.1800530716488682
$ws.Range("P11").Value = 0.1912195981346403
$ws.Range("Q11").Value = 62.85398969148444
$ws.Range("R11").Value = 565.68590722336
$ws.Range("S11").Value = 0.0414288037899113
$ws.Range("T11").Value = 0.04622652746355058

$ws.Range("G12").Value = 9.356602000000001
$ws.Range("H12").Value = 28.069806
$ws.Range("I12").Value = 0.2351065080674274
$ws.Range("J12").Value = 0.2470140528999518
$ws.Range("M12").Value = 7.207612333333334
$ws.Range("N12").Value = 21.622837
$ws.Range("O12").Value = 0.1890664495410015
$ws.Range("P12").Value = 0.2007919674510081
$ws.Range("Q12").Value = 67.43875997329134
$ws.Range("R12").Value = 606.948839759622
$ws.Range("S12").Value = 0.04445075274429132
$ws.Range("T12").Value = 0.04959843766982871

$ws.Range("G13").Value = 9.356602000000001
$ws.Range("H13").Value = 28.069806
$ws.Range("I13").Value = 0.2351065080674274
$ws.Range("J13").Value = 0.2470140528999518
$ws.Range("O13").Value = 0.2101379387118686
$ws.Range("P13").Value = 0.2231702676624555
$ws.Range("Q13").Value = 74.954821675004
$ws.Range("R13").Value = 674.593395075036
$ws.Range("S13").Value = 0.04940479698303448
$ws.Range("T13").Value = 0.05512619230207019

$ws.Range("G14").Value = 9.356602000000001
$ws.Range("H14").Value = 28.069806
$ws.Range("I14").Value = 0.2351065080674274
$ws.Range("J14").Value = 0.2470140528999518
$ws.Range("M14").Value = 9.361017666666667
$ws.Range("N14").Value = 28.083053
$ws.Range("O14").Value = 0.2455534915692039
$ws.Range("P14").Value = 0.2607822213107805
$ws.Range("Q14").Value = 87.58731662196868
$ws.Range("R14").Value = 788.285849597718
$ws.Range("S14").Value = 0.0577312239466
$ws.Range("T14").Value = 0.06441687341022809

$ws.Range("G15").Value = 9.356602000000001
$ws.Range("H15").Value = 28.069806
$ws.Range("I15").Value = 0.2351065080674274
$ws.Range("J15").Value = 0.2470140528999518
$ws.Range("M15").Value = 6.6785765
$ws.Range("N15").Value = 13.357153
$ws.Range("O15").Value = 0.1751890485290578
$ws.Range("P15").Value = 0.1240359454411156
$ws.Range("Q15").Value = 62.488782237053
$ws.Range("R15").Value = 374.932693422318
$ws.Range("S15").Value = 0.04118808545132184
$ws.Range("T15").Value = 0.03063862158868726

$ws.Range("G16").Value = 9.356602000000001
$ws.Range("H16").Value = 28.069806
$ws.Range("I16").Value = 0.2351065080674274
$ws.Range("J16").Value = 0.2470140528999518
$ws.Range("M16").Value = 6.864003333333333
$ws.Range("N16").Value = 20.59201
$ws.Range("O16").Value = 0.1800530716488682
$ws.Range("P16").Value = 0.1912195981346403
$ws.Range("Q16").Value = 64.22374731667333
$ws.Range("R16").Value = 578.0137258500599
$ws.Range("S16").Value = 0.04233164894217971
$ws.Range("T16").Value = 0.04723392792913758

$ws.Range("G17").Value = 5.7553975
$ws.Range("H17").Value = 11.510795
$ws.Range("I17").Value = 0.144617822663078
$ws.Range("J17").Value = 0.1012948976223954
$ws.Range("M17").Value = 7.207612333333334
$ws.Range("N17").Value = 21.622837
$ws.Range("O17").Value = 0.1890664495410015
$ws.Range("P17").Value = 0.2007919674510081
$ws.Range("Q17").Value = 41.48267400423583
$ws.Range("R17").Value = 248.896044025415
$ws.Range("S17").Value = 0.02734237827125835
$ws.Range("T17").Value = 0.02033920178634922

$ws.Range("G18").Value = 5.7553975
$ws.Range("H18").Value = 11.510795
$ws.Range("I18").Value = 0.144617822663078
$ws.Range("J18").Value = 0.1012948976223954
$ws.Range("O18").Value = 0.2101379387118686
$ws.Range("P18").Value = 0.2231702676624555
$ws.Range("Q18").Value = 46.105925343545
$ws.Range("R18").Value = 276.63555206127
$ws.Range("S18").Value = 0.03038969115541777
$ws.Range("T18").Value = 0.02260600941523101

$ws.Range("G19").Value = 5.7553975
$ws.Range("H19").Value = 11.510795
$ws.Range("I19").Value = 0.144617822663078
$ws.Range("J19").Value = 0.1012948976223954
$ws.Range("M19").Value = 9.361017666666667
$ws.Range("N19").Value = 28.083053
$ws.Range("O19").Value = 0.2455534915692039
$ws.Range("P19").Value = 0.2607822213107805
$ws.Range("Q19").Value = 53.87637767618917
$ws.Range("R19").Value = 323.258266057135
$ws.Range("S19").Value = 0.03551141129805476
$ws.Range("T19").Value = 0.02641590840941638

$ws.Range("G20").Value = 5.7553975
$ws.Range("H20").Value = 11.510795
$ws.Range("I20").Value = 0.144617822663078
$ws.Range("J20").Value = 0.1012948976223954
$ws.Range("M20").Value = 6.6785765
$ws.Range("N20").Value = 13.357153
$ws.Range("O20").Value = 0.1751890485290578
$ws.Range("P20").Value = 0.1240359454411156
$ws.Range("Q20").Value = 38.43786249165875
$ws.Range("R20").Value = 153.751449966635
$ws.Range("S20").Value = 0.02533545875268864
$ws.Range("T20").Value = 0.01256420839495482

$ws.Range("G21").Value = 5.7553975
$ws.Range("H21").Value = 11.510795
$ws.Range("I21").Value = 0.144617822663078
$ws.Range("J21").Value = 0.1012948976223954
$ws.Range("M21").Value = 6.864003333333333
$ws.Range("N21").Value = 20.59201
$ws.Range("O21").Value = 0.1800530716488682
$ws.Range("P21").Value = 0.1912195981346403
$ws.Range("Q21").Value = 39.50506762465833
$ws.Range("R21").Value = 237.03040574795
$ws.Range("S21").Value = 0.0260388831856585
$ws.Range("T21").Value = 0.01936956961644399

$ws.Range("G22").Value = 9.992936
$ws.Range("H22").Value = 29.978808
$ws.Range("I22").Value = 0.2510958880479564
$ws.Range("J22").Value = 0.2638132541845675
$ws.Range("M22").Value = 7.207612333333334
$ws.Range("N22").Value = 21.622837
$ws.Range("O22").Value = 0.1890664495410015
$ws.Range("P22").Value = 0.2007919674510081
$ws.Range("Q22").Value = 72.02520875981067
$ws.Range("R22").Value = 648.226878838296
$ws.Range("S22").Value = 0.04747380804757192
$ws.Range("T22").Value = 0.0529715823473722

$ws.Range("G23").Value = 9.992936
$ws.Range("H23").Value = 29.978808
$ws.Range("I23").Value = 0.2510958880479564
$ws.Range("J23").Value = 0.2638132541845675
$ws.Range("O23").Value = 0.2101379387118686
$ws.Range("P23").Value = 0.2231702676624555
$ws.Range("Q23").Value = 80.05243098827199
$ws.Range("R23").Value = 720.4718788944481
$ws.Range("S23").Value = 0.05276477233342367
$ws.Range("T23").Value = 0.05887527454927334

$ws.Range("G24").Value = 9.992936
$ws.Range("H24").Value = 29.978808
$ws.Range("I24").Value = 0.2510958880479564
$ws.Range("J24").Value = 0.2638132541845675
$ws.Range("M24").Value = 9.361017666666667
$ws.Range("N24").Value = 28.083053
$ws.Range("O24").Value = 0.2455534915692039
$ws.Range("P24").Value = 0.2607822213107805
$ws.Range("Q24").Value = 93.54405043786934
$ws.Range("R24").Value = 841.896453940824
$ws.Range("S24").Value = 0.06165747202884564
$ws.Range("T24").Value = 0.06879780643747709

$ws.Range("G25").Value = 9.992936
$ws.Range("H25").Value = 29.978808
$ws.Range("I25").Value = 0.2510958880479564
$ws.Range("J25").Value = 0.2638132541845675
$ws.Range("M25").Value = 6.6785765
$ws.Range("N25").Value = 13.357153
$ws.Range("O25").Value = 0.1751890485290578
$ws.Range("P25").Value = 0.1240359454411156
$ws.Range("Q25").Value = 66.738587535604
$ws.Range("R25").Value = 400.431525213624
$ws.Range("S25").Value = 0.04398924971668029
$ws.Range("T25").Value = 0.03272232640268017

$ws.Range("G26").Value = 9.992936
$ws.Range("H26").Value = 29.978808
$ws.Range("I26").Value = 0.2510958880479564
$ws.Range("J26").Value = 0.2638132541845675
$ws.Range("M26").Value = 6.864003333333333
$ws.Range("N26").Value = 20.59201
$ws.Range("O26").Value = 0.1800530716488682
$ws.Range("P26").Value = 0.1912195981346403
$ws.Range("Q26").Value = 68.59154601378667
$ws.Range("R26").Value = 617.3239141240799
$ws.Range("S26").Value = 0.04521058592143488
$ws.Range("T26").Value = 0.05044626444776473
